# Ajuste no algoritmo para considerar o preço de fechamento
# Updates simulation results (probabilities, balances, best-balance dates)
# in the "SimulacaoPeloDia" sheet to reflect the new algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    # Force text storage so numeric-/date-looking strings (e.g. "0.9992",
    # "2025-04-04") are not auto-converted by Excel into numbers/dates.
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Row 2 - PETRE301
Set-TextValue "H2" "0.9992"
Set-TextValue "I2" "0.5883"
Set-TextValue "N2" "R$ -375.97"
Set-TextValue "O2" "R$ -17.93"
Set-TextValue "P2" "2025-04-04"

# Row 3 - PETRE306
Set-TextValue "H3" "0.9976"
Set-TextValue "I3" "0.4645"
Set-TextValue "N3" "R$ -330.90"
Set-TextValue "O3" "R$ -71.70"
Set-TextValue "P3" "2025-04-16"

# Row 4 - PETRE312
Set-TextValue "H4" "0.9952"
Set-TextValue "I4" "0.3847"
Set-TextValue "N4" "R$ -366.46"
Set-TextValue "O4" "R$ -77.41"

# Row 5 - PETRE316
Set-TextValue "H5" "0.9911"
Set-TextValue "I5" "0.3107"
Set-TextValue "J5" "ITM → ATM"
Set-TextValue "N5" "R$ -363.35"
Set-TextValue "O5" "R$ -69.38"
Set-TextValue "P5" "2025-04-02"

# Row 6 - PETRE321
Set-TextValue "H6" "0.9843"
Set-TextValue "I6" "0.2447"
$ws.Range("L6").Value = 30
Set-TextValue "N6" "R$ -372.29"
Set-TextValue "O6" "R$ -108.90"

# Row 7 - PETRF303
Set-TextValue "H7" "0.6628"
Set-TextValue "N7" "R$ 929.08"
Set-TextValue "O7" "R$ 929.08"
Set-TextValue "P7" "2025-06-18"

# Row 8 - PETRF321
Set-TextValue "H8" "0.4271"
Set-TextValue "N8" "R$ 237.74"
Set-TextValue "O8" "R$ 908.51"
Set-TextValue "P8" "2025-06-10"

# Row 9 - PETRF326
Set-TextValue "H9" "0.3690"
Set-TextValue "I9" "0.9982"
Set-TextValue "N9" "R$ 172.17"
Set-TextValue "O9" "R$ 853.99"
Set-TextValue "P9" "2025-06-10"

# Row 10 - PETRF331
Set-TextValue "H10" "0.5522"
Set-TextValue "N10" "R$ 1735.24"
Set-TextValue "O10" "R$ 1735.24"
Set-TextValue "P10" "2025-06-18"

# Row 11 - PETRF342
Set-TextValue "H11" "0.2569"
Set-TextValue "I11" "0.7320"
Set-TextValue "N11" "R$ 242.57"
Set-TextValue "O11" "R$ 550.33"
Set-TextValue "P11" "2025-06-10"

# Row 12 - PETRF376
Set-TextValue "H12" "0.2016"
Set-TextValue "I12" "0.1700"
Set-TextValue "N12" "R$ 254.26"
Set-TextValue "O12" "R$ 353.57"
Set-TextValue "P12" "2025-06-10"
